$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per the diff
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 9
$ws.Range("A3").Value = 3
$ws.Range("B5").Value = 1

# Remove row 6 entirely (shifts dimension from A1:B6 to A1:B5)
$ws.Rows.Item(6).Delete()
